$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.863.20'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '1.894.53'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7938'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3208'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07116'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08068'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7752'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.43%  '

$ws.Range("D13").Value = '1.919.89'
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.329'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").Value = '29.899.78'
$ws.Range("E16").Value = '  +0.09%  '

$ws.Range("E17").Value = '  +0.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.938'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007763'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("D21").Value = '2.165.35'
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.154'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +19.08%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1625'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.322'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.075'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.45%  '

$ws.Range("E31").Value = '  +1.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.479'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05666'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.111'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.270'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7401'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.96%  '

$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.696'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01936'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.781'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4463'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.878'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8466'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.891'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.33%  '

$ws.Range("D47").Value = '1.030.56'
$ws.Range("E47").Value = '  +5.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.938'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.502'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("D51").Value = '2.060.10'
$ws.Range("E51").Value = '  +0.19%  '
